# RPA datasets push 2024-04-23
# The "확정공모가" (confirmed offering price) column (D) previously held
# "-" (not yet decided) for SK증권스팩12호 (row 10) and 민테크 (row 14).
# Both IPOs have now priced: fill in the confirmed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Assign as text (apostrophe-prefix keeps "2000"/"10500" from being
# auto-coerced to numbers), then restore the default "Normal" style so
# the cells keep looking exactly like their unstyled neighbours.
$ws.Range("D10").Value = "'2000"
$ws.Range("D10").Style = "Normal"

$ws.Range("D14").Value = "'10500"
$ws.Range("D14").Style = "Normal"
